$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "02030" sheet after the existing "06020" sheet.
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "02030"

# Activating the new sheet makes it the workbook's active tab (activeTab=1)
# and flips tabSelected from sheet1 to sheet2 automatically.
$ws2.Activate()

# ---------------------------------------------------------------------------
# 2. Column widths (bestFit-style widths copied from the source workbook).
#    ColumnWidth differs from the stored OOXML "width" by a constant 5/6
#    padding offset in this engine, so we pre-subtract it.
# ---------------------------------------------------------------------------
$offset = 5.0 / 6.0
$ws2.Columns.Item(1).ColumnWidth = 25.6640625 - $offset
$ws2.Columns.Item(2).ColumnWidth = 28        - $offset
$ws2.Columns.Item(3).ColumnWidth = 22.6640625 - $offset
$ws2.Columns.Item(4).ColumnWidth = 7.83203125 - $offset
$ws2.Columns.Item(5).ColumnWidth = 15.5       - $offset
$ws2.Columns.Item(6).ColumnWidth = 51        - $offset
$ws2.Columns.Item(7).ColumnWidth = 6.1640625 - $offset

# ---------------------------------------------------------------------------
# 3. Row heights.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(1).RowHeight = 16
$ws2.Rows.Item(2).RowHeight = 16
$ws2.Rows.Item(3).RowHeight = 16

# ---------------------------------------------------------------------------
# 4. Cell values - header row (row 1), reusing the same shared-string labels
#    as sheet "06020".
# ---------------------------------------------------------------------------
$ws2.Cells.Item(1,1).Value = "omschrijving"
$ws2.Cells.Item(1,2).Value = "inhoud"
$ws2.Cells.Item(1,3).Value = "weergave"
$ws2.Cells.Item(1,4).Value = "uitlijnen"
$ws2.Cells.Item(1,5).Value = "regel verwijderen"
$ws2.Cells.Item(1,6).Value = "regel template"
$ws2.Cells.Item(1,7).Value = "PTEST"

# Row 2
$ws2.Cells.Item(2,1).Value = "Verzekerd bedrag accessoires"
$ws2.Cells.Item(2,2).Value = 11239
$ws2.Cells.Item(2,3).Value = ""
$ws2.Cells.Item(2,4).Value = ""
$ws2.Cells.Item(2,5).Value = ""
$ws2.Cells.Item(2,6).Value = "09 Verzekerd bedrag accessoires € 11239"
$ws2.Cells.Item(2,7).Value = "x"

# Row 3
$ws2.Cells.Item(3,1).Value = ""
$ws2.Cells.Item(3,2).Value = 11239
$ws2.Cells.Item(3,3).Value = "Getal inclusief decimalen"
$ws2.Cells.Item(3,4).Value = "Rechts"
$ws2.Cells.Item(3,5).Value = "verwijderen"
$ws2.Cells.Item(3,6).Value = ""
$ws2.Cells.Item(3,7).Value = ""

# ---------------------------------------------------------------------------
# 5. Formatting - size 12 / black font, "@" (text) number format on every
#    used cell of the new sheet.
# ---------------------------------------------------------------------------
$used = $ws2.Range("A1:G3")
$used.NumberFormat = "@"
$used.Font.Size = 12
$used.Font.Color = 0

# ---------------------------------------------------------------------------
# 6. Selection on the new sheet.
# ---------------------------------------------------------------------------
$null = $ws2.Range("E6").Select()

# ---------------------------------------------------------------------------
# 7. Keep sheet1's own selection untouched (still G2) - nothing to do, it
#    only loses tabSelected which Activate() on ws2 already handled.
# ---------------------------------------------------------------------------
